$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.905.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.664.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.662.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("E10").Value = "  -4.03%  "

$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("E13").Value = "  -2.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.171.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("E15").Value = "  -4.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "71.858.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.673.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.88%  "

$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("E26").Value = "  -2.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.805.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0971"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.30%  "

$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("E39").Value = "  -3.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.47%  "

$ws.Range("E41").Value = "  -4.87%  "

$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("E51").Value = "  -1.96%  "
